$d = $word.ActiveDocument

# --- Step 1: insert a new paragraph right after "10-15 minute, MAXIM..."
#     containing the "10-15 minute speech..." text (this paragraph used to
#     exist two paragraphs later; a brand-new one is inserted here instead).
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()
$d.Paragraphs.Item(3).Range.Text = "10-15 minute speech ambele situatii – de invatat un speech de vreo 12 minute, estimativ"

# --- Step 2: the paragraph that used to hold "10-15 minute speech..." (now
#     shifted down to index 5, right after the blank paragraph) gets new text.
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "Unele aplicatii – nu toate (slide 3)"

# --- Step 3: insert three brand-new paragraphs right after it.
$p5.Range.InsertParagraphAfter()
$d.Paragraphs.Item(6).Range.Text = "Foarte pe scurt de ce am folosit alea – stocare date, securizare, algoritmi"

$d.Paragraphs.Item(6).Range.InsertParagraphAfter()
$d.Paragraphs.Item(7).Range.Text = "Baze de date, doar spun ca exista pentru conturi/profiluri, prin sql/nosql"

$d.Paragraphs.Item(7).Range.InsertParagraphAfter()
$d.Paragraphs.Item(8).Range.Text = "FARA PE PARTEA. Se poate remarca un exemplu de pagina pentru inregistrare"

# --- Step 4: remove the old trio of paragraphs (blank, "Cum ar trebui
#     prezentat?", and the numbered "Primele 2 slide-uri..." item) which are
#     now at indices 9-11, replacing them with a single new paragraph.
$p9 = $d.Paragraphs.Item(9)
$p11 = $d.Paragraphs.Item(11)
$oldRange = $d.Range($p9.Range.Start, $p11.Range.End)
$oldRange.Delete()

$p8 = $d.Paragraphs.Item(8)
$p8.Range.InsertParagraphAfter()
$newP = $d.Paragraphs.Item(9)
$newP.Range.Text = "Graficul din figura pune in evidenta..."
$newP.Range.LanguageID = "en-US"
